$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "https://www.youtube.com/watch?v=hY9nnU4PTFw&index=9&list=PLf0swTFhTI8rJvGpOp-LujOcpk-Rlz-yE"

$ws.Range("B12").Select()
